$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.682.80'
$ws.Range('D3').Value = '2.508.46'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('D5').Value = '596.47'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').Value = '176.42'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.63%  '
$ws.Range('D9').Value = '2.506.90'
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  +13.69%  '
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('D14').Value = '2.966.95'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('E16').Value = '  +4.14%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '69.453.48'
$ws.Range('E17').Value = '  +2.56%  '
$ws.Range('D18').Value = '2.493.04'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '363.73'
$ws.Range('E19').Value = '  +3.77%  '
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = '7.54'
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '70.70'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').Value = '9.04'
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('D27').Value = '1.67'
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('D28').Value = '2.640.41'
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').Value = '512.07'
$ws.Range('E30').Value = '  +1.42%  '
$ws.Range('D31').Value = '0.0₃0895'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('E33').Value = '  -1.76%  '
$ws.Range('D34').Value = '1.78'
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.119'
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '161.61'
$ws.Range('E37').Value = '  -0.83%  '
$ws.Range('D38').Value = '18.80'
$ws.Range('E38').Value = '  +2.66%  '
$ws.Range('D39').Value = '18.90'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('D43').Value = '4.82'
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('E44').Value = '  -2.47%  '
$ws.Range('D45').Value = '38.87'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('E46').Value = '  -3.63%  '
$ws.Range('D47').Value = '150.06'
$ws.Range('E47').Value = '  +3.71%  '
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('E51').Value = '  -0.46%  '
